# Trade #19 closed at 2026-02-17 13:18:03 - unknown UNKNOWN +0.000%
#
# Updates the Summary / Strategy Status roll-up numbers for the new
# MarketMaking trade, and appends the new trade (#19, row 20) to the
# "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Summary sheet
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.25               # Current Capital
$summary.Range("B4").Value = -0.76                 # Total P&L $
$summary.Range("B5").Value = -0.8                  # Total P&L %
$summary.Range("B6").Value = 19                    # Total Trades
$summary.Range("B8").Value = 12                    # Losing Trades
$summary.Range("B9").Value = 31.58                 # Win Rate %

# ------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.25                  # Capital
$status.Range("D4").Value = 19                     # Trades
$status.Range("E4").Value = -0.76                  # P&L $
$status.Range("F4").Value = -0.75                  # P&L %
$status.Range("G4").Value = 31.58                  # Win Rate %

# ------------------------------------------------------------------
# Append the new trade row (row 20, trade #19) to a trades sheet.
# Helper writes all 17 columns (A..Q); the Date column (B) needs to be
# forced to text first so the "yyyy-mm-dd" literal isn't auto-converted
# into a date serial number, then the style is reset back to the
# workbook default ("Normal") so no stray number-format sticks to the
# cell - matching how the rest of the sheet's text cells look.
# ------------------------------------------------------------------
function Add-TradeRow($ws) {
    $ws.Range("A20").Value = 19
    $ws.Range("B20").NumberFormat = "@"
    $ws.Range("B20").Value = "2026-02-17"
    $ws.Range("B20").Style = "Normal"
    $ws.Range("C20").Value = "13:17:57"
    $ws.Range("D20").Value = "MarketMaking"
    $ws.Range("E20").Value = "DOWN"
    $ws.Range("F20").Value = 0.777033
    $ws.Range("G20").Value = 0.75
    $ws.Range("H20").Value = "CLOSED"
    $ws.Range("I20").Value = -3.4789
    $ws.Range("J20").Value = -0.03
    $ws.Range("K20").Value = 99.25
    $ws.Range("L20").Value = 0
    $ws.Range("M20").Value = 0
    $ws.Range("N20").Value = 0.6
    $ws.Range("O20").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P20").Value = "early_exit"
    $ws.Range("Q20").Value = 0.13
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking
